$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "MuSCs" sending-cluster block (rows 6:9) — the data now only
# covers the single Wnt5b -> Fzd2 pair sent from FAPs.
$ws.Rows("6:9").Delete()

# Row 2 (FAPs -> Wnt5b/Fzd2 -> ECs): refreshed TPM-derived values
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("M2").Value = 0.1811433333333334
$ws.Range("N2").Value = 0.5434300000000001
$ws.Range("O2").Value = 0.0111261749556462
$ws.Range("P2").Value = 0.01112617495564619
$ws.Range("Q2").Value = 0.04070882434888889
$ws.Range("R2").Value = 0.3663794191400001
$ws.Range("S2").Value = 0.0111261749556462
$ws.Range("T2").Value = 0.01112617495564619

# Row 3 (FAPs -> Wnt5b/Fzd2 -> FAPs): refreshed TPM-derived values
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("O3").Value = 0.8246098959508241
$ws.Range("P3").Value = 0.8246098959508241
$ws.Range("Q3").Value = 3.017110511423555
$ws.Range("R3").Value = 27.15399460281199
$ws.Range("S3").Value = 0.8246098959508241
$ws.Range("T3").Value = 0.8246098959508241

# Row 4 (FAPs -> Wnt5b/Fzd2 -> MuSCs): refreshed TPM-derived values
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("M4").Value = 2.659118666666667
$ws.Range("N4").Value = 7.977356
$ws.Range("O4").Value = 0.1633282272592126
$ws.Range("P4").Value = 0.1633282272592126
$ws.Range("Q4").Value = 0.5975908289431111
$ws.Range("R4").Value = 5.378317460488
$ws.Range("S4").Value = 0.1633282272592126
$ws.Range("T4").Value = 0.1633282272592126

# Row 5 (FAPs -> Wnt5b/Fzd2 -> Resolving-Mac): refreshed TPM-derived values
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("M5").Value = 0.015234
$ws.Range("N5").Value = 0.045702
$ws.Range("O5").Value = 0.0009357018343171013
$ws.Range("P5").Value = 0.0009357018343171013
$ws.Range("Q5").Value = 0.003423577444
$ws.Range("R5").Value = 0.030812196996
$ws.Range("S5").Value = 0.0009357018343171013
$ws.Range("T5").Value = 0.0009357018343171013
